# Adds rows 119-133 to Sheet1 (case 21TRD09437), adding "Dismissed" finding
# entries and related duplicate case rows, per the JailCC dialog change.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rowsData = @(
    @{ Row = 119; A = '21TRD09437'; B = 'Hemmeter'; C = 'DUS'; D = '4510.11'; E = 'M1'; F = 'Guilty'; G = 'Guilty'; HasHI = $true }
    @{ Row = 120; A = '21TRD09437'; B = 'Hemmeter'; C = '1ST SPEED 1 YR SCHOOL >35MPHM4'; D = '4511.21B1A'; E = 'M4'; F = 'Dismissed'; HasHI = $true }
    @{ Row = 121; A = '21TRD09437'; B = 'Hemmeter'; C = 'RECKLESS OPERATION 1ST IN 1 YR'; D = '4511.20'; E = 'MM'; F = 'Guilty'; G = 'Guilty'; HasHI = $true }
    @{ Row = 122; A = '21TRD09437'; B = 'Hemmeter'; C = 'DUS'; D = '4510.11'; E = 'M1'; HasHI = $false }
    @{ Row = 123; A = '21TRD09437'; B = 'Hemmeter'; C = '1ST SPEED 1 YR SCHOOL >35MPHM4'; D = '4511.21B1A'; E = 'M4'; HasHI = $false }
    @{ Row = 124; A = '21TRD09437'; B = 'Hemmeter'; C = 'RECKLESS OPERATION 1ST IN 1 YR'; D = '4511.20'; E = 'MM'; HasHI = $false }
    @{ Row = 125; A = '21TRD09437'; B = 'Bunner'; C = 'DUS'; D = '4510.11'; E = 'M1'; HasHI = $false }
    @{ Row = 126; A = '21TRD09437'; B = 'Bunner'; C = '1ST SPEED 1 YR SCHOOL >35MPHM4'; D = '4511.21B1A'; E = 'M4'; HasHI = $false }
    @{ Row = 127; A = '21TRD09437'; B = 'Bunner'; C = 'RECKLESS OPERATION 1ST IN 1 YR'; D = '4511.20'; E = 'MM'; HasHI = $false }
    @{ Row = 128; A = '21TRD09437'; B = 'Hemmeter'; C = 'DUS'; D = '4510.11'; E = 'M1'; F = 'Guilty'; G = 'Guilty'; J = 'None'; K = 'None'; HasHI = $true }
    @{ Row = 129; A = '21TRD09437'; B = 'Hemmeter'; C = '1ST SPEED 1 YR SCHOOL >35MPHM4'; D = '4511.21B1A'; E = 'M4'; F = 'Dismissed'; J = 'None'; K = 'None'; HasHI = $true }
    @{ Row = 130; A = '21TRD09437'; B = 'Hemmeter'; C = 'RECKLESS OPERATION 1ST IN 1 YR'; D = '4511.20'; E = 'MM'; F = 'Guilty'; G = 'Guilty'; J = 'None'; K = 'None'; HasHI = $true }
    @{ Row = 131; A = '21TRD09437'; B = 'Hemmeter'; C = 'DUS'; D = '4510.11'; E = 'M1'; F = 'No Contest'; G = 'Guilty'; HasHI = $true }
    @{ Row = 132; A = '21TRD09437'; B = 'Hemmeter'; C = '1ST SPEED 1 YR SCHOOL >35MPHM4'; D = '4511.21B1A'; E = 'M4'; F = 'Dismissed'; HasHI = $true }
    @{ Row = 133; A = '21TRD09437'; B = 'Hemmeter'; C = 'RECKLESS OPERATION 1ST IN 1 YR'; D = '4511.20'; E = 'MM'; F = 'No Contest'; G = 'Guilty'; HasHI = $true }
)

$colIndex = @{ A=1; B=2; C=3; D=4; E=5; F=6; G=7; H=8; I=9; J=10; K=11 }

foreach ($entry in $rowsData) {
    $r = $entry.Row
    foreach ($col in @("A","B","C","D","E","F","G","J","K")) {
        if ($entry.ContainsKey($col)) {
            $cell = $ws.Cells.Item($r, $colIndex[$col])
            $cell.NumberFormat = "@"
            $cell.Value = $entry[$col]
        }
    }
    if ($entry.HasHI) {
        $ws.Cells.Item($r, $colIndex["H"]).Value = 0
        $hCell = $ws.Cells.Item($r, $colIndex["I"])
        $hCell.NumberFormat = "@"
        $hCell.Value = "0"
    }
}

# Row 132 column G is an empty (but present) text cell in the source data
$g132 = $ws.Cells.Item(132, $colIndex["G"])
$g132.NumberFormat = "@"